$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 74; this shifts the existing rows
# 74..136 down to 75..137 (and the sheet's dimension grows to A1:T137).
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new record.
$ws.Range("A74").Value = 4
$ws.Range("B74").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C74").Value = "Los Lagos"
$ws.Range("D74").Value = 44447
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100102
$ws.Range("H74").Value = "Cítricos"
$ws.Range("I74").Value = 100102006
$ws.Range("J74").Value = "Pomelo"
$ws.Range("K74").Value = "Start Ruby"
$ws.Range("L74").Value = "Primera"
$ws.Range("M74").Value = 50
$ws.Range("N74").Value = 12000
$ws.Range("O74").Value = 12000
$ws.Range("P74").Value = 12000
$ws.Range("Q74").Value = '$/caja 14 kilos empedrada'
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 857
$ws.Range("T74").Value = 14
